$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - No behaviours
$ws.Range("C2").Value = 88.63929303514631
$ws.Range("D2").Value = 40.03950428415587

# Row 3 - B1 (moving)
$ws.Range("B3").Value = 12.5
$ws.Range("C3").Value = 85.08306737582647
$ws.Range("D3").Value = 76.33610488677604

# Row 4 - B2 (requesting)
$ws.Range("B4").Value = 9.090909090909092
$ws.Range("C4").Value = 88.00492926060745
$ws.Range("D4").Value = 71.1588080641254

# Row 5 - B3 (notifying)
$ws.Range("C5").Value = 89.94321738516827
$ws.Range("D5").Value = 40.29332510135731

# Row 6 - B1 (moving) and B2 (requesting)
$ws.Range("B6").Value = 14.28571428571429
$ws.Range("C6").Value = 88.15907410538675
$ws.Range("D6").Value = 88.15907410538675

# Row 7 - B1 (moving) and B3 (notifying)
$ws.Range("B7").Value = 12.5
$ws.Range("C7").Value = 87.07803508584608
$ws.Range("D7").Value = 78.95681347173027

# Row 8 - B2 (requesting) and B3 (notifying)
$ws.Range("B8").Value = 9.090909090909092
$ws.Range("C8").Value = 89.53267275816702
$ws.Range("D8").Value = 74.28823824637688

# Row 9 - All behaviours
$ws.Range("B9").Value = 16.66666666666667
$ws.Range("C9").Value = 86.33706528481468
$ws.Range("D9").Value = 90.27271577374452
